# início das queries de análise serial
#
# Fixes the ordering of the bottom (tied/zero-value) ranking rows 23-26
# for UFs RO / SE / AC / TO across the "avg-arrecad", "max-arrecad" and
# "tx-sucesso" sheets. "tot-arrecad" keeps its current order (RO, SE, AC, TO).

$wb = $excel.ActiveWorkbook

# --- avg-arrecad: RO, SE, AC, TO  ->  AC, RO, SE, TO (row 26 "TO" stays put)
$wsAvg = $wb.Worksheets.Item("avg-arrecad")
$wsAvg.Range("A23").Value2 = "AC"
$wsAvg.Range("A24").Value2 = "RO"
$wsAvg.Range("A25").Value2 = "SE"

# --- max-arrecad: swap rows 25/26 (TO, AC) -> (AC, TO)
$wsMax = $wb.Worksheets.Item("max-arrecad")
$wsMax.Range("A25").Value2 = "AC"
$wsMax.Range("A26").Value2 = "TO"

# --- tx-sucesso: swap rows 25/26 (AC, TO) -> (TO, AC); B column (0) unaffected
$wsTx = $wb.Worksheets.Item("tx-sucesso")
$wsTx.Range("A25").Value2 = "TO"
$wsTx.Range("A26").Value2 = "AC"
